# Apply "Cập nhật plan mới" changes to the Project Plan and Gantt sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Project Plan and Gantt")

# Row 16: Responsible changes from "Lê Khánh" to "Trần Minh Trung"
$ws.Range("B16").Value = "Trần Minh Trung"

# Row 18: Start/End dates shift (9/16/2016-9/25/2016 -> 9/20/2016-9/26/2016)
$ws.Range("C18").Value = 42633
$ws.Range("D18").Value = 42639

# Row 20: Responsible changes from "Trần Minh Trung" to "Lê Khánh", End date shifts
$ws.Range("B20").Value = "Lê Khánh"
$ws.Range("D20").Value = 42640

# Row 21: Responsible changes from "Lê Khánh" to "Trần Minh Trung"
$ws.Range("B21").Value = "Trần Minh Trung"

# Update the sheet view: scroll position, zoom, and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$ws.Range("B26").Select()
$excel.ActiveWindow.Zoom = 100

$wb.Save()
